# RLS include Manager Role to view all
# Adds a "Manager" column to the User sheet and three new users
# (user1/user2/user3 @majesco.com), two of which are flagged as Managers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User")

# --- New B5 user (row 4: UserID 4, non-manager) -----------------------
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "user1@majesco.com"
$ws.Range("C5").Value = 0
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:user1@majesco.com")
$ws.Range("B5").Style = "Hyperlink"

# --- New B6/B7 users (rows 5 & 6: UserID 5 & 6, managers) --------------
$ws.Range("A6").Value = 5
$ws.Range("C6").Value = 1
$ws.Range("A7").Value = 6
$ws.Range("C7").Value = 1

# Stale multi-cell hyperlink left over from copy/pasting the B6 hyperlink
# down onto B7 (matches the leftover range-hyperlink entry in the sheet).
$ws.Hyperlinks.Add($ws.Range("B6:B7"), "mailto:user2@majesco.com", [Type]::Missing, [Type]::Missing, "user1@majesco.com")

# Real per-cell values/hyperlinks for the two new manager rows.
$ws.Range("B6").Value = "user2@majesco.com"
$ws.Range("B7").Value = "user3@majesco.com"
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:user2@majesco.com")
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:user3@majesco.com")
$ws.Range("B6").Style = "Hyperlink"
$ws.Range("B7").Style = "Hyperlink"

# --- Manager column for existing rows + header -------------------------
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("C1").Value = "Manager"

# --- Selection bookkeeping (Project sheet moves to A5, User stays active
#     and ends up selecting C3) -----------------------------------------
$ws2 = $wb.Worksheets.Item("Project")
[void]$ws2.Activate()
[void]$ws2.Range("A5").Select()

[void]$ws.Activate()
[void]$ws.Range("C3").Select()
